$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# Column A ("IA Control") cells: the comma-separated list of controls was
# re-ordered (content/items unchanged, only their order within the cell).
# -----------------------------------------------------------------------
$ws.Range("A2").Value = 'AU-4,AU-4 (1)'
$ws.Range("A3").Value = 'SC-5,SC-5 (2),CM-6 b'
$ws.Range("A4").Value = 'AU-8 b,AC-6 (9),CM-5 (1),AU-7 b,AC-6 (8),AU-12 (3),AU-7 a'
$ws.Range("A5").Value = 'CM-7 b,AC-17 (1),AC-17 (9),CM-6 b'
$ws.Range("A8").Value = 'IA-2 (11),IA-2 (12)'
$ws.Range("A10").Value = 'CM-7 (2),CM-7 (5) (b)'
$ws.Range("A15").Value = 'IA-2,AU-3 (1),IA-8'
$ws.Range("A17").Value = 'MA-4 (1) (a),AU-3 (1),AU-3,AU-12 c,AU-12 a'
$ws.Range("A19").Value = 'IA-5 (1) (b),CM-6 b,IA-5 (1) (a)'
$ws.Range("A21").Value = 'SC-10,MA-4 e,MA-4 (7),AC-12'
$ws.Range("A22").Value = 'AU-14 (1),MA-4 (1) (a),AU-3 (1),AU-7 (1),AU-3,CM-5 (1),AU-6 (4),CM-6 b,AU-7 a,AU-12 a'
$ws.Range("A25").Value = 'MA-4 (1) (a),AU-3 (1),AU-3,AU-12 c,AU-12 a'
$ws.Range("A31").Value = 'MA-4 (1) (a),AU-3 (1),AC-2 (4),AU-3,AU-12 c,AU-12 a'
$ws.Range("A34").Value = 'AC-11 a,AC-11 b'
$ws.Range("A38").Value = 'SI-11 b,AU-9'
$ws.Range("A45").Value = 'AC-8 c 1, AC-8 c 2, AC-8 c 3,AC-8 a,AC-8 b'
$ws.Range("A50").Value = 'IA-2 (5),CM-6 b'
$ws.Range("A56").Value = 'AU-12 c,MA-4 (1) (a)'
$ws.Range("A63").Value = 'AU-5 a,AU-5 (1)'
$ws.Range("A65").Value = 'IA-2 (2),CM-6 b'
$ws.Range("A67").Value = 'MA-4 (1) (a),AU-3 (1),AU-3,AU-12 c,AU-12 a'
$ws.Range("A69").Value = 'AU-8 b,AU-12 c,CM-5 (1),AU-7 b,CM-6 b,AU-12 (3),AU-7 a,AU-12 a'
$ws.Range("A77").Value = 'MA-4 (1) (a),AU-3 (1),AC-2 (4),AU-3,AU-12 c,AU-12 a'
$ws.Range("A79").Value = 'AU-9 (3),AU-9'
$ws.Range("A80").Value = 'IA-2 (4),IA-2 (2),IA-2 (3),IA-2 (1)'
$ws.Range("A86").Value = 'MA-4 (1) (a),AU-3 (1),AU-3,AU-12 c,AU-12 a'
$ws.Range("A88").Value = 'AC-2 (4),AU-12 c,AC-6 (9),CM-5 (1)'
$ws.Range("A89").Value = 'IA-2 (4),IA-2 (2),IA-2 (3),IA-2,IA-2 (5)'
$ws.Range("A90").Value = 'IA-2 (11),IA-2 (12)'
$ws.Range("A91").Value = 'AU-9 (3),AU-9'
$ws.Range("A97").Value = 'AU-8 b,AU-8 (1) (b),AU-8 (1) (a)'
$ws.Range("A101").Value = 'AC-3 (4),IA-11'
$ws.Range("A102").Value = 'MA-4 (1) (a),AU-3 (1),AU-3,AU-12 c,AU-12 a'
$ws.Range("A111").Value = 'AU-5 a,AU-5 b'
$ws.Range("A119").Value = 'MA-4 (1) (a),AU-3 (1),AU-3,AU-12 c,AU-12 a'
$ws.Range("A124").Value = 'MA-4 (1) (a),AU-3 (1),AU-3,AU-12 c,AU-12 a'
$ws.Range("A128").Value = 'IA-5 (1) (c),CM-7 a,CM-6 b'
$ws.Range("A139").Value = 'SI-6 d,SI-6 b,CM-3 (5)'
$ws.Range("A148").Value = 'AU-14 (1),MA-4 (1) (a),AU-3 (1),AU-3,AU-12 c,AU-12 a'
$ws.Range("A157").Value = 'MA-4 (1) (a),AU-3 (1),AU-3,AU-12 c,AU-12 a'
$ws.Range("A175").Value = 'CM-7 a,SI-16'

# -----------------------------------------------------------------------
# Row 16 (CtrlAltDelBurstAction): fill in the previously empty Fix (M) cell
# -----------------------------------------------------------------------
$ws.Range("M16").Value = 'Edit "/etc/systemd/system.conf" and add or edit the following line:
CtrlAltDelBurstAction=none'

# -----------------------------------------------------------------------
# Row 43 (maxlogins): update Check (K) text and fill in Fix (M) cell
# -----------------------------------------------------------------------
$ws.Range("K43").Value = 'Run the following command to ensure the "maxlogins" value is
configured for all users on the system:
 # grep "maxlogins" /etc/security/limits.conf /etc/security/limits.d/*.conf 
You should receive output similar to the following:
 *\t\thard\tmaxlogins\t10 

If maxlogins is not equal to or less than the expected value then this is a finding.'

$ws.Range("M43").Value = 'Configure the operating system to limit the number of concurrent sessions to "10" for all accounts and/or account types.
Add the following line to the top of the /etc/security/limits.conf or in a ".conf" file defined in /etc/security/limits.d/ :
* hard maxlogins 10'

# -----------------------------------------------------------------------
# Row 97 (maxpoll / NTP): update Check (K) text and fill in Fix (M) cell
# -----------------------------------------------------------------------
$ws.Range("K97").Value = 'Verify Red Hat Enterprise Linux 9 is securely comparing internal information system clocks at a regular interval with an NTP server with the following commands:
To verify that "maxpoll" has been set properly, perform the following:
 $ sudo grep maxpoll /etc/ntp.conf /etc/chrony.conf 
The output should return:
 maxpoll 16 .

If maxpoll does not exist or maxpoll has not been set to the expected value of 16 then this is a finding.'

$ws.Range("M97").Value = 'Configure Red Hat Enterprise Linux 9 to securely compare internal information system clocks at a regular interval with an NTP server by adding/modifying the following line in the /etc/chrony.conf file.
server [ntp.server.name] iburst maxpoll 16'

# -----------------------------------------------------------------------
# Row 177 (SSSD offline credentials expiration): fill in Fix (M) cell
# -----------------------------------------------------------------------
$ws.Range("M177").Value = 'Configure the SSSD to prohibit the use of cached authentications after one day.
Add or change the following line in "/etc/sssd/sssd.conf" just below the line "[pam]".
offline_credentials_expiration = 1'
